# Add a new "Driver License" mapping row to the IEPD mapping sheet
# (Business Class / Business Attribute / Data Container-Type / NIEM Element /
#  NIEM Element Type / IEP Path), matching the other sections already
# present on the sheet. This inserts a brand-new section between the
# existing "Contact Information" (row 39) and "Registrant Residence
# Location" (old row 41) sections, pushing everything below down by two
# rows (one new content row + the usual blank separator row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two blank rows at row 41: one for the new content, one to keep the
# blank-row separator pattern used between every section on this sheet.
$ws.Range("A41:A42").EntireRow.Insert()

# Populate the new row. Columns, in sheet order: A=Business Class,
# B=Business Attribute, C=Data Container/Type, D=NIEM Element,
# E=NIEM Element Type, F=IEP Path.
# NOTE: values are written in A, C, B, D, E, F order so that new shared
# strings land in the same index order as the target workbook.
$ws.Range("A41").Value = "Driver License"
$ws.Range("C41").Value = "nc:IdentificationType"
$ws.Range("B41").Value = "Driver License ID"
$ws.Range("D41").Value = "nc:IdentificationID"
$ws.Range("E41").Value = "niem-xsd:string"
$ws.Range("F41").Value = "exchange:FirearmRegistrationQueryResults/nc:DriverLicense/nc:DriverLicenseIdentification/nc:IdentificationID"

# Leave the cursor where the author left it.
$null = $ws.Range("F46").Select()
